# Apply the commit-metadata and prose edits described by the diff.

$d = $word.ActiveDocument

# 1. Commit ID
$d.Content.Find.Execute("b2b401e", $true, $false, $false, $false, $false, $true, 1, $false, "79f515b", 2) | Out-Null

# 2. Commit date
$d.Content.Find.Execute("2023-07-21", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-01", 2) | Out-Null

# 3. FirstParagraph: trim the "In fact... never heard of the term. I thought..." sentence
$old3 = "But that isn" + [char]8217 + "t how I have always viewed them. In fact, when I began my PhD I had never heard of the term. I thought of reporting guidelines as similar sets of recommendations that authors should adhere to. I didn" + [char]8217 + "t think about how guidance"
$new3 = "But that isn" + [char]8217 + "t how I have always viewed them. When I began my PhD, I didn" + [char]8217 + "t think about how guidance"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 4. Append sentence after "authors used guidance."
$old4 = "authors used guidance."
$new4 = "authors used guidance. I merely thought of reporting guidelines as similar sets of recommendations that authors should adhere to."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# 5. Rewrite opening of the "Having identified..." paragraph
$old5 = "Having identified so many influences that may be limiting the success of the current system, I started looking for a framework that could help me work out how this system could be changed to address the barriers I identified. I considered the MRC guidance"
$new5 = "I began looking for a framework that could help me understand and improve this complex system. I considered the MRC guidance"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null
